$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Date" property value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-01-15T07:47:50+00:00"

# --- Elements sheet: update row 7 (Extension.value[x]:valueAddress) ---
$elements = $wb.Worksheets.Item("Elements")

# Definition now mirrors the parent Extension.value[x] row's definition
$elements.Range("M7").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."

# Comments cleared
$elements.Range("N7").Value = ""

# Condition(s) cleared
$elements.Range("AI7").Value = ""

# Mapping: RIM Mapping now "N/A"
$elements.Range("AK7").Value = "N/A"
